$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows: Day index (A), B, C, D, E, F for rows 9 through 56.
# (Rows 2-8 are unchanged by this edit.)
$updates = @(
    @(9, 0, 0, 34, 13, 21),
    @(10, 0, 0, 40, 19, 21),
    @(11, 0, 0, 45, 24, 21),
    @(12, 0, 0, 53, 32, 21),
    @(13, 0, 21, 54, 33, 0),
    @(14, 0, 21, 56, 26, 9),
    @(15, 0, 21, 63, 25, 17),
    @(16, 0, 21, 67, 24, 22),
    @(17, 0, 21, 79, 27, 31),
    @(18, 0, 21, 99, 45, 33),
    @(19, 0, 30, 114, 58, 26),
    @(20, 0, 38, 120, 59, 23),
    @(21, 0, 43, 139, 74, 22),
    @(22, 0, 52, 153, 80, 21),
    @(23, 1, 54, 177, 83, 40),
    @(24, 1, 56, 215, 107, 52),
    @(25, 1, 61, 254, 138, 55),
    @(26, 1, 65, 303, 170, 68),
    @(27, 1, 73, 351, 206, 72),
    @(28, 1, 94, 417, 247, 76),
    @(29, 4, 108, 482, 286, 88),
    @(30, 6, 116, 576, 342, 118),
    @(31, 6, 133, 669, 397, 139),
    @(32, 6, 145, 751, 421, 185),
    @(33, 6, 170, 872, 478, 224),
    @(34, 6, 196, 1036, 592, 248),
    @(35, 6, 234, 1195, 665, 296),
    @(36, 6, 272, 1381, 754, 355),
    @(37, 6, 330, 1610, 900, 380),
    @(38, 6, 394, 1799, 976, 429),
    @(39, 7, 444, 2054, 1096, 514),
    @(40, 9, 530, 2356, 1241, 585),
    @(41, 13, 627, 2749, 1450, 672),
    @(42, 13, 710, 3161, 1660, 791),
    @(43, 14, 823, 3646, 1948, 875),
    @(44, 17, 958, 4159, 2219, 982),
    @(45, 21, 1115, 4729, 2516, 1098),
    @(46, 22, 1299, 5391, 2816, 1276),
    @(47, 27, 1501, 6153, 3175, 1477),
    @(48, 35, 1698, 7054, 3651, 1705),
    @(49, 38, 1940, 8037, 4127, 1970),
    @(50, 45, 2213, 9170, 4739, 2218),
    @(51, 52, 2575, 10397, 5319, 2503),
    @(52, 56, 2978, 11751, 5937, 2836),
    @(53, 64, 3403, 13228, 6623, 3202),
    @(54, 76, 3910, 14773, 7214, 3649),
    @(55, 98, 4431, 16530, 7902, 4197),
    @(56, 116, 5078, 18420, 8625, 4717)
)

foreach ($row in $updates) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove the now-obsolete trailing rows (formerly rows 57-59, days 55-57),
# shrinking the sheet's populated range down to row 56.
$ws.Range("A57:F59").Delete()
